$wb = $excel.ActiveWorkbook

$oldText = "January 30 2026 16.19.47 EST"
$newText = "February 02 2026 12.49.33 EST"

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            $val = $cell.Value2
            if ($val -ne $null -and $val -is [string] -and $val.Contains($oldText)) {
                $cell.Value2 = $val.Replace($oldText, $newText)
            }
        }
    }
}
